# Update "paises.xlsx" style Covid-19 country table + timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 08:22"

# 2) Updated per-country figures (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes).
$updates = @{
    "Austria"  = @(13561, 1,   6064, 7178, 261, 0, 319)
    "Hungria"  = @(1310,  120, 115,  1110, 17,  8, 85)
    "Lituania" = @(1026,  27,  54,   949,  14,  1, 23)
    "Kuwait"   = @(993,   0,   133,  859,  26,  0, 1)
    "Taiwan"   = @(385,   3,   99,   280,  0,   0, 6)
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 4; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($updates.ContainsKey($name)) {
        $vals = $updates[$name]
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
        $ws.Cells.Item($r, 5).Value = $vals[3]
        $ws.Cells.Item($r, 6).Value = $vals[4]
        $ws.Cells.Item($r, 7).Value = $vals[5]
        $ws.Cells.Item($r, 8).Value = $vals[6]
    }
}

# 3) Re-sort the country table by "Casos totales" (column B) descending, as the
#    source data is always kept in that order; updated countries move to their
#    new rank.
$sortRange = $ws.Range("A4:H" + $lastRow)
$keyRange = $ws.Range("B4:B" + $lastRow)
$sortRange.Sort($keyRange, 2)
